$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.354.00"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "1.881.16"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'0.7116"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").Value = "'242.78"
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "'0.08026"
$ws.Range("E8").Value = "  +3.19%  "
$ws.Range("D9").Value = "'0.3143"
$ws.Range("E9").Value = "  +0.99%  "
$ws.Range("E10").Value = "  -0.11%  "
$ws.Range("D11").Value = "'0.08328"
$ws.Range("E11").Value = "  -1.50%  "
$ws.Range("D12").Value = "1.896.76"
$ws.Range("E12").Value = "  +1.07%  "
$ws.Range("D13").Value = "'5.266"
$ws.Range("E13").Value = "  +0.64%  "
$ws.Range("D14").Value = "'94.80"
$ws.Range("E14").Value = "  +4.03%  "
$ws.Range("D15").Value = "'0.7185"
$ws.Range("E15").Value = "  +0.92%  "
$ws.Range("D16").Value = "'6.376"
$ws.Range("E16").Value = "  +5.44%  "
$ws.Range("D17").Value = "'0.000008673"
$ws.Range("E17").Value = "  +5.57%  "
$ws.Range("D18").Value = "29.371.32"
$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("D19").Value = "'243.16"
$ws.Range("E19").Value = "  +0.90%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "'13.35"
$ws.Range("E20").Value = "  +0.77%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.148.10"
$ws.Range("E21").Value = "  +1.29%  "
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D23").Value = "'7.830"
$ws.Range("E23").Value = "  +0.74%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("E25").Value = "  -1.91%  "
$ws.Range("D26").Value = "'163.45"
$ws.Range("D27").Value = "'9.091"
$ws.Range("E27").Value = "  +0.41%  "
$ws.Range("E28").Value = "  +0.56%  "
$ws.Range("D29").Value = "'1.509"
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("D30").Value = "'4.434"
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("D31").Value = "'4.357"
$ws.Range("E31").Value = "  +1.32%  "
$ws.Range("E32").Value = "  -6.44%  "
$ws.Range("D33").Value = "'0.05396"
$ws.Range("E33").Value = "  +2.28%  "
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("D35").Value = "'0.7767"
$ws.Range("E35").Value = "  +4.28%  "
$ws.Range("D36").Value = "'1.183"
$ws.Range("E36").Value = "  +0.42%  "
$ws.Range("D37").Value = "'2.689"
$ws.Range("E37").Value = "  -0.35%  "
$ws.Range("D38").Value = "'0.01889"
$ws.Range("E38").Value = "  +0.96%  "
$ws.Range("D39").Value = "1.272.42"
$ws.Range("E39").Value = "  +4.31%  "
$ws.Range("D40").Value = "'2.746"
$ws.Range("E40").Value = "  +0.78%  "
$ws.Range("D41").Value = "'6.536"
$ws.Range("E41").Value = "  +0.88%  "
$ws.Range("D42").Value = "'0.9192"
$ws.Range("E42").Value = "  +3.41%  "
$ws.Range("D43").Value = "'113.51"
$ws.Range("E43").Value = "  +4.32%  "
$ws.Range("D44").Value = "'74.51"
$ws.Range("E44").Value = "  +2.52%  "
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "'0.00000000128"
$ws.Range("E46").Value = "  +4.76%  "
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "2.043.37"
$ws.Range("E47").Value = "  +1.27%  "
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("D49").Value = "'0.5222"
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("D50").Value = "'9.591"
$ws.Range("E50").Value = "  +2.43%  "
$ws.Range("E51").Value = "  +1.41%  "
